$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-parsed as numbers (losing formatting / introducing float error).
$textCells = @("D5", "D6", "D7", "D8", "D10", "D13", "D15", "D19", "D20", "D21", "D24", "D25", "D28", "D29", "D32", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values scraped for this run.
$ws.Range("D2").Value = '59.017.47'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '2.514.87'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '530.68'
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").Value = '138.99'
$ws.Range("E6").Value = '  -2.97%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.563'
$ws.Range("E8").Value = '  -1.53%  '
$ws.Range("D9").Value = '2.518.13'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = '0.100'
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("D13").Value = '0.356'
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").Value = '2.959.47'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '23.00'
$ws.Range("E15").Value = '  -2.10%  '
$ws.Range("D16").Value = '58.972.28'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").Value = '2.506.82'
$ws.Range("D19").Value = '11.01'
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("D20").Value = '4.25'
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").Value = '322.16'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").Value = '62.26'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").Value = '0.424'
$ws.Range("E25").Value = '  -3.70%  '
$ws.Range("E26").Value = '  +2.02%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").Value = '7.80'
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").Value = '6.81'
$ws.Range("E29").Value = '  +2.22%  '
$ws.Range("D30").Value = '0.0₃0770'
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("E31").Value = '  -1.76%  '
$ws.Range("D32").Value = '162.27'
$ws.Range("E32").Value = '  +3.84%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").Value = '1.12'
$ws.Range("E34").Value = '  -5.96%  '
$ws.Range("E35").Value = '  -0.64%  '
$ws.Range("D36").Value = '18.46'
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").Value = '4.21'
$ws.Range("E37").Value = '  -3.02%  '
$ws.Range("D38").Value = '1.57'
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").Value = '37.04'
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").Value = '3.64'
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("D41").Value = '0.803'
$ws.Range("E41").Value = '  -1.92%  '
$ws.Range("D42").Value = '5.20'
$ws.Range("E42").Value = '  -8.64%  '
$ws.Range("D43").Value = '279.55'
$ws.Range("E43").Value = '  -5.44%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '10.85'
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").Value = '0.595'
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").Value = '0.0931'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("D48").Value = '121.66'
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("D49").Value = '18.33'
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").Value = '0.0510'
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("E51").Value = '  -2.39%  '
